$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows before row 946, shifting existing data (old 946-956) down to 954-964
for ($i = 0; $i -lt 8; $i++) {
    $ws.Rows.Item(946).EntireRow.Insert()
}

$values = @(
    @{Row=946; D=44628; K='Moscatel rosada'; L='Primera'; M=250; N=16000; O=17000; P=16500; Q='$/bandeja 18 kilos'; R='Provincia de Limarí'; S=917; T=18}
    @{Row=947; D=44628; K='Red Globe'; L='Segunda'; M=200; N=7000; O=7000; P=7000; Q='$/bandeja 10 kilos'; R='Provincia de Limarí'; S=700; T=10}
    @{Row=948; D=44628; K='Rosada pastilla'; L='Primera'; M=300; N=10000; O=10000; P=10000; Q='$/bandeja 10 kilos'; R='Provincia de Limarí'; S=1000; T=10}
    @{Row=949; D=44628; K='Thompson seedless'; L='Primera'; M=300; N=11000; O=12000; P=11500; Q='$/bandeja 18 kilos'; R='Provincia de Limarí'; S=639; T=18}
    @{Row=950; D=44628; K='Thompson seedless'; L='Primera'; M=600; N=11000; O=12000; P=11500; Q='$/bandeja 18 kilos'; R='Provincia de San Felipe de Aconcagua'; S=639; T=18}
    @{Row=951; D=44628; K='Thompson seedless'; L='Primera'; M=1000; N=11000; O=12000; P=11500; Q='$/bandeja 18 kilos'; R='Región de O''Higgins'; S=639; T=18}
    @{Row=952; D=44628; K='Thompson seedless'; L='Segunda'; M=250; N=9000; O=10000; P=9500; Q='$/bandeja 18 kilos'; R='Provincia de Limarí'; S=528; T=18}
    @{Row=953; D=44628; K='Thompson seedless'; L='Segunda'; M=800; N=9000; O=10000; P=9500; Q='$/bandeja 18 kilos'; R='Región de O''Higgins'; S=528; T=18}
)

foreach ($row in $values) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100109
    $ws.Cells.Item($r, 8).Value = "Uva"
    $ws.Cells.Item($r, 9).Value = 100109001
    $ws.Cells.Item($r, 10).Value = "Uva"
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}
